# Auto-generated PowerShell COM-interop script
# Updates column F ('想去人数' / want-to-go count) values across all 4 sheets
# to match the commit '456a3b4' generated-output refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 334
$ws.Range("F4").Value = 1313
$ws.Range("F5").Value = 376
$ws.Range("F6").Value = 356
$ws.Range("F7").Value = 3905
$ws.Range("F9").Value = 771
$ws.Range("F10").Value = 2316
$ws.Range("F11").Value = 352
$ws.Range("F12").Value = 227
$ws.Range("F14").Value = 198
$ws.Range("F15").Value = 187
$ws.Range("F16").Value = 2243
$ws.Range("F17").Value = 322
$ws.Range("F18").Value = 26
$ws.Range("F20").Value = 345
$ws.Range("F21").Value = 234
$ws.Range("F22").Value = 44
$ws.Range("F23").Value = 276

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 26
$ws.Range("F5").Value = 36
$ws.Range("F7").Value = 131
$ws.Range("F12").Value = 5
$ws.Range("F22").Value = 62

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6399
$ws.Range("F4").Value = 2114
$ws.Range("F5").Value = 340
$ws.Range("F6").Value = 12

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6399
$ws.Range("F4").Value = 2114
$ws.Range("F5").Value = 340
$ws.Range("F10").Value = 334
$ws.Range("F11").Value = 1313
$ws.Range("F12").Value = 376
$ws.Range("F13").Value = 26
$ws.Range("F14").Value = 36
$ws.Range("F16").Value = 12
$ws.Range("F17").Value = 356
$ws.Range("F18").Value = 3905
$ws.Range("F19").Value = 131
$ws.Range("F24").Value = 771
$ws.Range("F25").Value = 2316
$ws.Range("F26").Value = 352
$ws.Range("F28").Value = 227
$ws.Range("F30").Value = 198
$ws.Range("F31").Value = 187
$ws.Range("F32").Value = 5
$ws.Range("F34").Value = 2243
$ws.Range("F35").Value = 322
$ws.Range("F38").Value = 26
$ws.Range("F40").Value = 345
$ws.Range("F41").Value = 234
$ws.Range("F42").Value = 44
$ws.Range("F49").Value = 62
$ws.Range("F50").Value = 276
